$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Create the "Normal 2" cell style (mirrors the xfId=1 / "Normal 2"
# cellStyle that appears when data is pasted in from another workbook).
# ------------------------------------------------------------------
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.ColorIndex = 8

# ------------------------------------------------------------------
# Add the new worksheet after the last existing sheet.
# ------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Duke 18-19"

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Range("A1:B1").ColumnWidth = 35.166666666666664
$ws.Range("C1:J1").ColumnWidth = 9.877604166666666
$ws.Range("K1:P1").ColumnWidth = 4.675781
$ws.Range("Q1").ColumnWidth = 245.30729166666666

# ------------------------------------------------------------------
# Title / header block, rows 1-7
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Platform Report 1 (R4)"
$ws.Range("B1").Value = "Total Searches, Result Clicks and Record Views by Month and Platform"
$ws.Range("A2").Value = "FLORIDA STATE UNIV"
$ws.Range("A3").Value = " "
$ws.Range("A4").Value = "Period covered by Report:"
$ws.Range("A5").Value = "2018-07-01 to 2019-06-30"
$ws.Range("A6").Value = "Date run:"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2019-09-05"

$ws.Rows.Item(1).RowHeight = 15.95
$ws.Rows.Item(2).RowHeight = 15.95
$ws.Rows.Item(3).RowHeight = 15.95
$ws.Rows.Item(4).RowHeight = 15.95
$ws.Rows.Item(5).RowHeight = 15.95
$ws.Rows.Item(6).RowHeight = 15.95
$ws.Rows.Item(7).RowHeight = 12
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15

$titleRange = $ws.Range("A1:B7")
$titleRange.Style = "Normal 2"
$titleRange.Font.Name = "Arial"
$titleRange.Font.Size = 9
$titleRange.Font.Bold = $true
$titleRange.Font.ColorIndex = 8
$titleRange.HorizontalAlignment = -4131
$titleRange.VerticalAlignment = -4160

# ------------------------------------------------------------------
# Header row 8
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Platform"
$ws.Range("B8").Value = "Publisher"
$ws.Range("C8").Value = "User Activity"
$ws.Range("D8").Value = "Reporting Period Total"
$ws.Range("E8").Value = "Jul-2018"
$ws.Range("F8").Value = "Aug-2018"
$ws.Range("G8").Value = "Sep-2018"
$ws.Range("H8").Value = "Oct-2018"
$ws.Range("I8").Value = "Nov-2018"
$ws.Range("J8").Value = "Dec-2018"
$ws.Range("K8").Value = "Jan-2019"
$ws.Range("L8").Value = "Feb-2019"
$ws.Range("M8").Value = "Mar-2019"
$ws.Range("N8").Value = "Apr-2019"
$ws.Range("O8").Value = "May-2019"
$ws.Range("P8").Value = "Jun-2019"

$headerRange = $ws.Range("A8:P8")
$headerRange.Style = "Normal 2"
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9
$headerRange.Font.Bold = $true
$headerRange.Font.ColorIndex = 8
$headerRange.Interior.Color = 11439959
$headerRange.HorizontalAlignment = -4131
$ws.Range("E8").WrapText = $true
$ws.Range("K8:P8").WrapText = $true

# ------------------------------------------------------------------
# Data rows 9-12
# ------------------------------------------------------------------
$ws.Range("A9:A12").Value = "Silverchair"
$ws.Range("B9:B12").Value = "Duke University Press"
$ws.Range("C9").Value = "Regular Searches"
$ws.Range("C10").Value = "Searches-federated and automated"
$ws.Range("C11").Value = "Result Clicks"
$ws.Range("C12").Value = "Record Views"

$row9 = @(176, 10, 12, 14, 26, 24, 0, 14, 14, 30, 18, 8, 6)
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $row9[$i]
}
for ($r = 10; $r -le 12; $r++) {
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value = 0
    }
}

$dataRange = $ws.Range("A9:P12")
$dataRange.Style = "Normal 2"
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 9
$dataRange.Font.ColorIndex = 8
$dataRange.Interior.Color = 14602940
$dataRange.HorizontalAlignment = -4131
$dataRange.VerticalAlignment = -4160
$dataRange.WrapText = $true

# ------------------------------------------------------------------
# View / page setup
# ------------------------------------------------------------------
$ws.Range("D33").Select() | Out-Null
$ws.PageSetup.PrintGridlines = $true
$ws.PageSetup.Orientation = 2

$ws.Activate()

Write-Host "done"
